$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.25"
$ws.Range("E2").Value = "'-0.12%"
$ws.Range("D3").Value = "'44.10"
$ws.Range("E3").Value = "'-0.59%"
$ws.Range("D4").Value = "'5.496"
$ws.Range("E4").Value = "'-1.64%"
$ws.Range("D5").Value = "'0.08043"
$ws.Range("E5").Value = "'-0.53%"
$ws.Range("D6").Value = "'2.061"
$ws.Range("E6").Value = "'3.69%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9582"
$ws.Range("E7").Value = "'0.42%"
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D8").Value = "'0.1126"
$ws.Range("E8").Value = "'-3.99%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1881"
$ws.Range("E9").Value = "'1.58%"
$ws.Range("B10").Value = "MCDex"
$ws.Range("C10").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D10").Value = "'10.17"
$ws.Range("E10").Value = "'-0.36%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09941"
$ws.Range("E11").Value = "'1.52%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04683"
$ws.Range("E12").Value = "'3.15%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1056"
$ws.Range("E13").Value = "'-1.14%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001268"
$ws.Range("E14").Value = "'-0.98%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04106"
$ws.Range("E15").Value = "'-2.28%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006095"
$ws.Range("E16").Value = "'3.81%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.340"
$ws.Range("E17").Value = "'-1.11%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.432"
$ws.Range("E18").Value = "'2.86%"
$ws.Range("D19").Value = "'2.636"
$ws.Range("E19").Value = "'2.78%"
$ws.Range("D20").Value = "'0.3320"
$ws.Range("E20").Value = "'-4.62%"
$ws.Range("D21").Value = "'0.1384"
$ws.Range("E21").Value = "'-2.46%"
$ws.Range("D22").Value = "'0.2579"
$ws.Range("E22").Value = "'2.92%"
$ws.Range("D23").Value = "'0.001313"
$ws.Range("E23").Value = "'5.33%"
$ws.Range("D24").Value = "'0.004375"
$ws.Range("E24").Value = "'0.45%"
$ws.Range("D25").Value = "'0.0001288"
$ws.Range("E25").Value = "'8.05%"
$ws.Range("D26").Value = "'0.0003747"
$ws.Range("E26").Value = "'-5.82%"
$ws.Range("D38").Value = "'0.02656"
$ws.Range("E38").Value = "'-1.15%"
$ws.Range("D39").Value = "'0.05608"
$ws.Range("E39").Value = "'0.74%"
$ws.Range("D40").Value = "'0.007645"
$ws.Range("E40").Value = "'0.63%"
$ws.Range("D41").Value = "'0.1413"
$ws.Range("E41").Value = "'0.22%"
$ws.Range("D42").Value = "'0.007386"
$ws.Range("E42").Value = "'-7.07%"
$ws.Range("D43").Value = "'0.001996"
$ws.Range("E43").Value = "'-1.13%"
$ws.Range("D44").Value = "'0.008700"
$ws.Range("E44").Value = "'3.51%"
$ws.Range("D45").Value = "'0.00007121"
$ws.Range("E45").Value = "'-0.80%"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.25%"
$ws.Range("D47").Value = "'0.0005814"
$ws.Range("E47").Value = "'0.04%"
$ws.Range("D48").Value = "'0.002525"
$ws.Range("E48").Value = "'11.31%"
$ws.Range("D49").Value = "'0.003447"
$ws.Range("E49").Value = "'-14.93%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.25%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.25%"
